$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of simulated "compras" quotation data
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 45413
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 45422
$ws.Range("E2").Value = "Solicitada"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 45413
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 45422
$ws.Range("E3").Value = "Aprovada"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 45413
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 45424
$ws.Range("E4").Value = "Solicitada"

$ws.Range("B2:B4").NumberFormat = "dd/mm/yy"
$ws.Range("D2:D4").NumberFormat = "dd/mm/yy"
